$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 1347.204755070514
$ws.Range("D8").Value = 24.50151034028677
